$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.405.60'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.646.38'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.01%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.11'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.59'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.17%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.645.40'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.99%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +8.04%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.51%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +2.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.24'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +2.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000192'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +2.30%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.127.30'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.344.29'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.636.90'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '365.32'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.51'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.70%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.90'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.53%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.32%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.17'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.31%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.97%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.97%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.17%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '575.67'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.04%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.56%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.09%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.61%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.74%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +6.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '159.44'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.31%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.35%  '
$ws.Range('B40').NumberFormat = '@'
$ws.Range('B40').Value = 'PolygonEcosystemToken'
$ws.Range('C40').NumberFormat = '@'
$ws.Range('C40').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.375'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +2.09%  '
$ws.Range('B41').NumberFormat = '@'
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').NumberFormat = '@'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.91'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.44'
$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'BabyDogeCoin'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0₆0338'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +10.48%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.68'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.45%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.73%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.33'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '157.40'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.56%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.77'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.73'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.19%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.88%  '
